# BBI-23-1.xlsx — "Add files via upload" edit
#
# Сheet "Лист1": three students' Лаба №1 grades (column B) that were
# marked with the placeholder "*" get updated:
#   - B6  (Воробьев Иван Александрович)      "*" -> "+*"  (typed with a
#   - B16 (Оганезов Михаил Алексеевич)       "*" -> "+*"   leading apostrophe,
#                                                           so Excel stores it
#                                                           as text with a
#                                                           quote-prefix style
#                                                           instead of trying
#                                                           to parse it as a
#                                                           formula)
#   - B20 (Стоценко Александр Сергеевич)     "*" -> 5     (an actual numeric
#                                                           grade was entered)
# B19 (Словогородская) keeps its "*" placeholder untouched.
# The active selection also moved on to B21 before the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Leading "'" mirrors typing "+*" directly into the cell: Excel would
# otherwise treat a leading "+" as the start of a formula, so it silently
# applies a quote-prefix (text) cell style instead.
$ws.Range("B6").Value = "'+*"
$ws.Range("B16").Value = "'+*"

# A real numeric grade replaces the placeholder asterisk.
$ws.Range("B20").Value = 5

# Move the active cell the way it was left when the sheet was saved.
$ws.Range("B21").Select()
